$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New consolidated text for A1 (previously the shared string content lived in A2)
$newText = @'
questions = [
    {
        "title": "The MEX number of a non-negative set of numbers is the smallest non-negative number that is not present in the set. For example, MEX({1 3 10})=0, and MEX({0 1 2 8})=3.Your task is to write the function updateMEX which takes a given array arr of length num and removes the minimum number of elements from it so that the MEX value of the modified array is not equal to the MEX value of the original array.The arguments are passed in the order: num, arrYour code should return the minimum number of elements that need to be removed from the array.If the task is not possible, then your code should return -1.Keep in mind:Array arr elements are non-negative integersArray elements are not necessarily distinct1 &lt= num &lt= 400 &lt= arr[i] &lt= 90Example 1Input:num = 4arr = 0 1 1 4Output:1 Explanation: The MEX of the input array is 2. If we remove the element \"0\" from it, we have a modified array of 1 1 4 with MEX = 0, which is not equal to 2. So, the answer is 1 as removing one element changed the MEX of the array.Example 2Input:num = 4arr = 1 2 3 4Output:-1Explanation: The MEX of the input array is 0. No matter what elements you remove from the array, the MEX remains unchanged. Hence we return -1.",
        "ques_type": null,
        "options": [],
        "score": null
    }
]
'@

# Remove row 2 entirely (shifts data up, drops the extra trailing row)
$ws.Range("A2").EntireRow.Delete()

# Clear the bold/border/centered formatting that used to be applied to A1
$ws.Range("A1").ClearFormats()

# Put the updated text into A1
$ws.Range("A1").Value = $newText

# Avoid pinning an explicit row height for the multi-line text
$ws.Rows.Item(1).AutoFit()
